# Add files via upload
#
# 1) Fix the typo in the "total accelleration - gravity" shared string
#    (used throughout column C of Sheet1) -> "total acceleration - gravity".
# 2) Update the selection/active-cell on Sheet1 to C4 (range C4:C82).
# 3) Remove Sheet2 entirely (its data was a duplicate listing of the
#    Description column and is no longer needed).

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$ws1 = $wb.Worksheets.Item("Sheet1")

# Fix the "accelleration" -> "acceleration" typo everywhere it occurs.
[void]$ws1.Cells.Replace("total accelleration - gravity", "total acceleration - gravity")

# Update the saved selection / active cell for Sheet1.
[void]$ws1.Range("C4:C82").Select()

# Remove Sheet2 from the workbook.
$ws2 = $wb.Worksheets.Item("Sheet2")
[void]$ws2.Delete()
